# suivi.xlsx - "modification v.1 du suivi"
#
# The Backlog sheet has four columns per task block: B=task name,
# C="A faire" (to do), D="En cours" (in progress), E="Fini" (done).
# Marking a task done means moving its label out of the "A faire" (C)
# column into the "Fini" (E) column (and the C cell goes back to being a
# plain empty cell).
#
# This edit marks four Jenkins-server tasks as completed:
#   row 32 - Afficher le contenu de /var/jenkins_home/secrets/initialAdminPassword
#   row 33 - Creer l'utilisateur userjob
#   row 34 - Lui donner le droit d'utiliser apt dans le fichier sudoers
#   row 36 - Ouvrir les ports de Jenkins et le port 22
# (row 35, "Installer un pare-feu", stays in the "A faire" column)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

$doneRows = 32,33,34,36

foreach ($r in $doneRows) {
    $todo = $ws.Cells.Item($r, 3)   # column C - "A faire"
    $done = $ws.Cells.Item($r, 5)   # column E - "Fini"

    # Carry the task label + its formatting (yellow fill/border) over to
    # the "Fini" column ...
    $todo.Copy($done)
    # ... and put the "A faire" cell back to a plain, empty cell.
    $todo.Clear()
}

$excel.CutCopyMode = $false

# Reflect where the user ended up looking / what was selected when the
# sheet was last saved.
$ws.Activate()
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B21").Select()
